$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet from "My Series" to "Data"
$ws.Name = "Data"

# 2. Update the label in A11
$ws.Range("A11").Value = "Function Information"

# 3. Insert 3 new rows before row 27 so the historical data series gains
#    three earlier observations (rows shift from 27:31 down to 30:34).
$ws.Rows("27:29").Insert()

# Copy number formats (date format for column A, decimal format for
# column B) down from the row that used to be row 27 (now row 30) so the
# freshly inserted rows keep the same look as the rest of the series.
$ws.Range("A30:B30").Copy()
$ws.Range("A27:B29").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 4. Populate the three newly inserted observations
$ws.Range("A27").Value = 37956
$ws.Range("B27").Value = 16421.36
$ws.Range("A28").Value = 38322
$ws.Range("B28").Value = 22352.63
$ws.Range("A29").Value = 38687
$ws.Range("B29").Value = 29205.22

# 5. Refresh the cached summary statistics (B17:B26) for the now 8-point
#    series (previously computed over the 5-point series).
$ws.Range("B17").Value = 68708.75249999999
$ws.Range("B18").Value = 3244857856.999478
$ws.Range("B19").Value = 56963.65382416649
$ws.Range("B20").Value = 1.182177702322342
$ws.Range("B21").Value = 0.6111834863010071
$ws.Range("B22").Value = 0.8290596430806468
$ws.Range("B23").Value = 16421.36
$ws.Range("B25").Value = 47666.46
$ws.Range("B26").Value = 8

# 6. Update the numeric display format used by the data column so values
#    show without a leading zero suppression change (0.000 -> ###0.000).
$ws.Range("B27:B34").NumberFormat = "###0.000"
